$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('Combined Ratio', 94.2, 'Intact', 'Q4 2024')
    ,@('Claims Ratio', 68.90000000000001, 'Intact', 'Q4 2024')
    ,@('Core Claim Ratio', 74.3, 'Intact', 'Q4 2024')
    ,@('CAT Loss Ratio', -0.5, 'Intact', 'Q4 2024')
    ,@('Expense Ratio', 25.3, 'Intact', 'Q4 2024')
    ,@('PYD Ratio', -4.9, 'Intact', 'Q4 2024')
    ,@('Gross Written Premium', 1575, 'Intact', 'Q4 2024')
    ,@('Underwriting Income', 97, 'Intact', 'Q4 2024')
    ,@('ROE', 14.2, 'Intact', 'Q4 2024')
    ,@('Combined Ratio', 96.09999999999999, 'Definity', 'Q4 2024')
    ,@('Combined Ratio', 95.40000000000001, 'Intact', '2024')
    ,@('Claims Ratio', 69.59999999999999, 'Intact', '2024')
    ,@('Core Claim Ratio', 71.8, 'Intact', '2024')
    ,@('CAT Loss Ratio', 1.6, 'Intact', '2024')
    ,@('Expense Ratio', 25.8, 'Intact', '2024')
    ,@('PYD Ratio', -3.8, 'Intact', '2024')
    ,@('Gross Written Premium', 6640, 'Intact', '2024')
    ,@('Gross Written Premium', 5956, 'Intact', '2023')
    ,@('Underwriting Income', 306, 'Intact', '2023')
    ,@('Underwriting Income', 292, 'Intact', '2024')
    ,@('Combined Ratio', 94.7, 'Intact', '2023')
    ,@('Claims Ratio', 68.90000000000001, 'Intact', '2023')
    ,@('Core Claim Ratio', 74.3, 'Intact', '2023')
    ,@('CAT Loss Ratio', 1.1, 'Intact', '2023')
    ,@('Expense Ratio', 25.8, 'Intact', '2023')
    ,@('PYD Ratio', -6.5, 'Intact', '2023')
    ,@('Combined Ratio', 95.2, 'Intact', 'Q4 2023')
    ,@('Claims Ratio', 69.8, 'Intact', 'Q4 2023')
    ,@('Core Claim Ratio', 75.8, 'Intact', 'Q4 2023')
    ,@('CAT Loss Ratio', -0.1, 'Intact', 'Q4 2023')
    ,@('Expense Ratio', 25.4, 'Intact', 'Q4 2023')
    ,@('PYD Ratio', -5.9, 'Intact', 'Q4 2023')
    ,@('Gross Written Premium', 1408, 'Intact', 'Q4 2023')
    ,@('Underwriting Income', 74, 'Intact', 'Q4 2023')
    ,@('ROE', 8.800000000000001, 'Intact', 'Q4 2023')
    ,@('Combined Ratio', 95.90000000000001, 'Definity', 'Q4 2023')
    ,@('Claims Ratio', 71.2, 'Definity', 'Q4 2023')
    ,@('Claims Ratio', 70.40000000000001, 'Definity', 'Q4 2024')
    ,@('Core Claim Ratio', 71.3, 'Definity', 'Q4 2024')
    ,@('Core Claim Ratio', 71.5, 'Definity', 'Q4 2023')
    ,@('CAT Loss Ratio', 1.1, 'Definity', 'Q4 2023')
    ,@('CAT Loss Ratio', 0.1, 'Definity', 'Q4 2024')
    ,@('Expense Ratio', 25.7, 'Definity', 'Q4 2024')
    ,@('Expense Ratio', 24.7, 'Definity', 'Q4 2023')
    ,@('PYD Ratio', -1.4, 'Definity', 'Q4 2023')
    ,@('PYD Ratio', -1, 'Definity', 'Q4 2024')
    ,@('Gross Written Premium', 438.7, 'Definity', 'Q4 2024')
    ,@('Gross Written Premium', 416, 'Definity', 'Q4 2023')
    ,@('Underwriting Income', 16.2, 'Definity', 'Q4 2023')
    ,@('Underwriting Income', 16.8, 'Definity', 'Q4 2024')
    ,@('ROE', 10.6, 'Definity', 'Q4 2024')
    ,@('ROE', 9.199999999999999, 'Definity', 'Q4 2023')
    ,@('Combined Ratio', 98.3, 'Definity', '2023')
    ,@('Combined Ratio', 96.7, 'Definity', '2024')
    ,@('Claims Ratio', 71.09999999999999, 'Definity', '2024')
    ,@('Claims Ratio', 71.8, 'Definity', '2023')
    ,@('Core Claim Ratio', 72.59999999999999, 'Definity', '2023')
    ,@('Core Claim Ratio', 70.90000000000001, 'Definity', '2024')
    ,@('CAT Loss Ratio', 1.1, 'Definity', '2024')
    ,@('CAT Loss Ratio', 0.9, 'Definity', '2023')
    ,@('Expense Ratio', 26.5, 'Definity', '2023')
    ,@('Expense Ratio', 25.6, 'Definity', '2024')
    ,@('PYD Ratio', -0.9, 'Definity', '2024')
    ,@('PYD Ratio', -1.7, 'Definity', '2023')
    ,@('Gross Written Premium', 1657.1, 'Definity', '2023')
    ,@('Gross Written Premium', 1867.4, 'Definity', '2024')
    ,@('Underwriting Income', 54.9, 'Definity', '2024')
    ,@('Underwriting Income', 26.5, 'Definity', '2023')
    ,@('ROE', 9.199999999999999, 'Definity', '2023')
    ,@('ROE', 10.6, 'Definity', '2024')
)

$startRow = 130
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]

    $dVal = $row[3]
    if ($dVal -match "^-?[0-9]+(\.[0-9]+)?$") {
        # Force pure-numeric-looking strings (e.g. "2024") to remain text
        $ws.Cells.Item($r, 4).NumberFormat = "@"
        $ws.Cells.Item($r, 4).Value = $dVal
        $ws.Cells.Item($r, 4).Style = "Normal"
    } else {
        $ws.Cells.Item($r, 4).Value = $dVal
    }
}
